$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: update maneuver-less config values, clear maneuver columns ---
$ws1.Range("B2").Value = 3000
$ws1.Range("C2").Value = 5000000
$ws1.Range("D2").Value = 12000
$ws1.Range("E2:H3").ClearContents()

# --- Sheet2: update maneuver-less config values, clear maneuver columns ---
$ws2.Range("B2").Value = 4000
$ws2.Range("C2").Value = 3000000
$ws2.Range("D2").Value = 13000
$ws2.Range("E2:H3").ClearContents()

# --- Selections ---
$ws1.Range("H18").Select() | Out-Null
$ws2.Range("E2:H4").Select() | Out-Null

# --- Activate Sheet2 (becomes the tabSelected / active tab) ---
$ws2.Activate()
